$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 12:46"

# Refresh province/city case counts (Casos totales, Casos activos, Recuperados, Muertes)
# and re-rank rows 4:63 by descending Casos totales, matching the updated source data.
$data = @(
    @("Madrid", 10575, 2063, 7249, 1263),
    @("Cataluña", 5925, 644, 5036, 245),
    @("Valencia/Valencia", 1084, 23, 1024, 37),
    @("Araba/Alava", 912, 21, 852, 60),
    @("Navarra", 886, 7, 855, 24),
    @("Bizkaia/Vizcaya", 860, 21, 836, 24),
    @("La Rioja", 747, 18, 707, 22),
    @("Alacant/Alicante", 603, 12, 541, 50),
    @("Asturias", 545, 12, 467, 10),
    @("A Coruña", 524, 19, 512, 12),
    @("Malaga", 505, 72, 464, 21),
    @("Ciudad Real", 505, 8, 457, 40),
    @("Toledo", 501, 22, 451, 28),
    @("Pontevedra", 445, 19, 442, 3),
    @("Albacete", 430, 8, 390, 32),
    @("Cantabria", 347, 11, 330, 6),
    @("Murcia", 345, 1, 213, 3),
    @("Granada", 335, 72, 276, 13),
    @("Zaragoza", 329, 0, 315, 14),
    @("Gipuzkoa/Guipuzcoa", 325, 21, 312, 13),
    @("Sevilla", 295, 72, 243, 4),
    @("Burgos", 269, 27, 175, 16),
    @("Salamanca", 265, 13, 180, 21),
    @("Guadalajara", 263, 2, 257, 4),
    @("Tenerife", 262, 4, 211, 11),
    @("Caceres", 243, 2, 231, 12),
    @("Valladolid", 241, 13, 193, 11),
    @("Castello/Castellon", 211, 1, 203, 7),
    @("Mallorca", 210, 10, 197, 9),
    @("Leon", 201, 3, 156, 14),
    @("Cordoba", 180, 72, 140, 3),
    @("Aragon", 174, 0, 163, 11),
    @("Jaen", 171, 72, 136, 3),
    @("Segovia", 157, 9, 140, 20),
    @("Ourense", 142, 19, 140, 2),
    @("Cadiz", 134, 72, 126, 2),
    @("Cuenca", 120, 8, 104, 8),
    @("Gran Canaria", 119, 0, 96, 11),
    @("Soria", 119, 5, 71, 8),
    @("Avila", 114, 14, 55, 9),
    @("Badajoz", 111, 5, 104, 2),
    @("Lugo", 62, 19, 60, 2),
    @("Zamora", 59, 3, 42, 3),
    @("Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3),
    @("Almeria", 58, 72, 53, 1),
    @("Teruel", 47, 0, 45, 2),
    @("Huelva", 47, 72, 37, 0),
    @("Palencia", 41, 2, 25, 0),
    @("Huesca", 37, 0, 37, 0),
    @("Melilla", 28, 0, 28, 0),
    @("Ibiza", 21, 10, 20, 1),
    @("Menorca", 15, 10, 13, 0),
    @("Fuerteventura", 14, 0, 13, 11),
    @("La Palma", 14, 0, 11, 11),
    @("Lanzarote", 9, 0, 4, 11),
    @("Arroyo de la Luz", 7, 0, 7, 0),
    @("Ceuta", 6, 0, 6, 0),
    @("La Gomera", 3, 2, 1, 11),
    @("El Hierro", 2, 0, 1, 11),
    @("Formentera", 0, 10, 0, 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 4 + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
